# Actualización automática 2025-11-14 08:30:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("E4").Value = 64.81999999999999
$ws1.Range("M4").Value = 1869.5

$ws1.Range("D24").Value = 457.92
$ws1.Range("E24").Value = 213.37
$ws1.Range("L24").Value = 2892.06

$ws1.Range("M29").Value = 128.3

$ws1.Range("L36").Value = 2311.1

# Row 56 "N de 54" summary counters shift right starting at column D
$ws1.Range("D56").Value = "6 de 54"
$ws1.Range("E56").Value = "2 de 54"
$ws1.Range("F56").Value = "0 de 54"
$ws1.Range("G56").Value = "0 de 54"
$ws1.Range("H56").Value = "1 de 54"
$ws1.Range("I56").Value = "3 de 54"
$ws1.Range("J56").Value = "0 de 54"
$ws1.Range("L56").Value = "6 de 54"
$ws1.Range("M56").Value = "12 de 54"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F4").Value = 3596.95
$ws2.Range("F24").Value = 3905.96
$ws2.Range("F29").Value = 3406.53
$ws2.Range("F36").Value = 2367.96
$ws2.Range("F60").Value = 24701.32

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Columns.Item(5).ColumnWidth = 22.1667

$ws3.Range("D3").Value = 3343.78
$ws3.Range("E3").Value = 3279.48
$ws3.Range("F3").Value = 0.504854105078164

$ws3.Range("D4").Value = 278.19
$ws3.Range("E4").Value = 566.51
$ws3.Range("F4").Value = 0.3293358588848112

$ws3.Range("D11").Value = 8742.02
$ws3.Range("E11").Value = 5493.969999999999
$ws3.Range("F11").Value = 0.6140788241632651

$ws3.Range("D12").Value = 9447.58
$ws3.Range("E12").Value = 55496.42
$ws3.Range("F12").Value = 0.1454727149544223

$ws3.Range("D14").Value = 24677.92
$ws3.Range("E14").Value = 74278.33685923838
$ws3.Range("F14").Value = 0.2493821086533561
